# Apply the "settings.xlsx" content/format cleanup described in the commit:
#   - rename the Simple_settings helper row to settings_content and give it
#     the new INSTALLED_APPS list (drop the old allauth/community entries,
#     add simple_history)
#   - apply uniform wrap-text formatting (Normal style) across the whole
#     A2:B22 settings table, which also drives Excel's row auto-fit heights
#   - move the active selection to A7:A8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------
$ws.Range("B1").Value = "settings_content"

$newInstalledApps = "[  'rest_framework',`n" + `
    "    'simple_history',`n" + `
    "    'django.contrib.admin',`n" + `
    "    'django.contrib.auth',`n" + `
    "    'django.contrib.contenttypes',`n" + `
    "    'django.contrib.sessions',`n" + `
    "    'django.contrib.messages',`n" + `
    "    'django.contrib.staticfiles',`n" + `
    "    'crispy_forms',`n" + `
    "    'corsheaders',`n" + `
    "    'debug_toolbar',`n" + `
    "    'django_filters',`n" + `
    "]"
$ws.Range("B7").Value = $newInstalledApps

# --- Formatting: reset to Normal style + wrap text on the whole table --
$table = $ws.Range("A2:B22")
$table.Style = "Normal"
$table.WrapText = $true

# --- Row heights (auto-fit result of the wrap-text formatting above) --
$ws.Rows.Item(2).RowHeight = 68
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 238
$ws.Rows.Item(8).RowHeight = 153
$ws.Rows.Item(9).RowHeight = 17
$ws.Rows.Item(10).RowHeight = 255
$ws.Rows.Item(11).RowHeight = 17
$ws.Rows.Item(12).RowHeight = 102
$ws.Rows.Item(13).RowHeight = 238
$ws.Rows.Item(14).RowHeight = 17
$ws.Rows.Item(15).RowHeight = 17
$ws.Rows.Item(16).RowHeight = 17
$ws.Rows.Item(17).RowHeight = 17
$ws.Rows.Item(18).RowHeight = 17
$ws.Rows.Item(19).RowHeight = 17
$ws.Rows.Item(20).RowHeight = 51
$ws.Rows.Item(21).RowHeight = 17
$ws.Rows.Item(22).RowHeight = 17

# --- Selection ---------------------------------------------------------
$ws.Range("A7:A8").Select()
